$wb = $excel.ActiveWorkbook

$wsTesting = $wb.Worksheets.Item("Testing & treatment")
$wsTesting.Rows("34:40").Delete()

$wsCascade = $wb.Worksheets.Item("Cascade")
$wsCascade.Rows("34:40").Delete()
